# Add the new "T1 - NIST sphere 1" .. "T1 - NIST sphere 14" columns
# (AD:AQ) to the right of the existing "T1 - cortical GM" column (AC),
# fix the casing of "t1map" -> "T1map" inside every NIFTI filename in
# column D, and extend the sheet's used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room for the 14 new columns. AD:AQ are currently past the
#    used range, so inserting there just appends the columns without
#    shifting any existing data.
$ws.Range("AD1:AQ1").EntireColumn.Insert()

# 2. Copy the header formatting (bold, centered, bordered) from the
#    last existing header cell (AC1) onto the new header cells so the
#    new header row matches the look of the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AQ1").PasteSpecial(-4122)

# 3. Fill in the new header labels.
$headers = @(
    "T1 - NIST sphere 1",
    "T1 - NIST sphere 2",
    "T1 - NIST sphere 3",
    "T1 - NIST sphere 4",
    "T1 - NIST sphere 5",
    "T1 - NIST sphere 6",
    "T1 - NIST sphere 7",
    "T1 - NIST sphere 8",
    "T1 - NIST sphere 9",
    "T1 - NIST sphere 10",
    "T1 - NIST sphere 11",
    "T1 - NIST sphere 12",
    "T1 - NIST sphere 13",
    "T1 - NIST sphere 14"
)

$columns = @("AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ")

for ($i = 0; $i -lt $columns.Length; $i++) {
    $ws.Range($columns[$i] + "1").Value = $headers[$i]
}

# 4. Fix the filename casing in column D for every data row
#    (t1map.nii.gz -> T1map.nii.gz).
$ws.Range("D2:D57").Replace("t1map.nii.gz", "T1map.nii.gz")
